# Building Single page applications.pptx -- "Adding some knockout action to flask"
#
# 1. Bump the cached datetimeFigureOut footer field from 9/2/2012 to 9/27/2012
#    on the slide master and every slide layout.
# 2. Append a new Title-and-Content slide (#8) asking how asp.net mvc helps
#    with ajax.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Refresh the "last modified/printed" date footer everywhere it is cached.
# ---------------------------------------------------------------------------
$newDate = "9/27/2012"

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name.StartsWith("Date Placeholder")) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$design = $p.Designs.Item(1)
$master = $design.SlideMaster
Update-DateShapes $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2. Add the new "How can asp.net mvc help with ajax?" slide at the end.
# ---------------------------------------------------------------------------
$newIndex = $p.Slides.Count + 1
$slide = $p.Slides.Add($newIndex, 2)

# Title: "How can asp.net mvc help with ajax?"
$titleRange = $slide.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "How can asp.net "
$run = $titleRange.InsertAfter("mvc")
$run = $run.InsertAfter(" help with ")
$run = $run.InsertAfter("ajax")
$run = $run.InsertAfter("?")

# Body bullets.
$bodyRange = $slide.Shapes.Item(2).TextFrame.TextRange
$bodyRange.Text = "Controllers can take and return JSON"
$run = $bodyRange.InsertAfter("`r")
$run = $run.InsertAfter("WebAPI")
$run = $run.InsertAfter(" allows for ")
$run = $run.InsertAfter("RESTful")
$run = $run.InsertAfter(" ")
$run = $run.InsertAfter("services")
$run = $run.InsertAfter("`r")
$run = $run.InsertAfter("JavaScript ")
$run = $run.InsertAfter("minification")
$run = $run.InsertAfter(" and concatenation ")
